$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Autonomous_temporary")

# --- 1. Fix typo in the JETSON_MS Mission_Select "Choices" text (row 7, col K) ---
$ws.Cells.Item(7, 11).Value = "0=Manual, 1=Acceleration, 2=Skidpadd, 3=Trackdrive, 4=Braketest, 5=Inspection, 6=Autocross"

# --- 2. Widen columns A and K by 1 character each ---
# ColumnWidth read back from COM is (stored XML width - 5/6), so subtract 5/6 from the
# desired stored width to land exactly on the target integer width.
$ws.Columns.Item(1).ColumnWidth = 22 - (5/6)
$ws.Columns.Item(11).ColumnWidth = 92 - (5/6)

# --- 3. Insert two new signal rows (EBS_Pressure_Front / EBS_Pressure_Rear) into the
#        RD_JETSON... wait, ACU_IGN message block, right after the "Emergency" signal (row 21) ---
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

# Copy the formatting (style) of the preceding data row onto the two new rows
$ws.Range("A21:K21").Copy()
$ws.Range("A22:K23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(22, 1).Value = "EBS_Pressure_Front"
$ws.Cells.Item(22, 2).Value = 24
$ws.Cells.Item(22, 3).Value = 8
$ws.Cells.Item(22, 4).Value = "Intel"
$ws.Cells.Item(22, 5).Value = $false
$ws.Cells.Item(22, 6).Value = 10
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 150
$ws.Cells.Item(22, 10).Value = "bar"

$ws.Cells.Item(23, 1).Value = "EBS_Pressure_Rear"
$ws.Cells.Item(23, 2).Value = 32
$ws.Cells.Item(23, 3).Value = 8
$ws.Cells.Item(23, 4).Value = "Intel"
$ws.Cells.Item(23, 5).Value = $false
$ws.Cells.Item(23, 6).Value = 10
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 150
$ws.Cells.Item(23, 10).Value = "bar"

# --- 4. Insert two new signal rows (Brake_pressure_Front / Brake_pressure_Rear) into the
#        VCU_HV message block, right after the "HV" signal (now row 35) ---
$ws.Rows.Item(36).Insert()
$ws.Rows.Item(36).Insert()

$ws.Range("A35:K35").Copy()
$ws.Range("A36:K37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(36, 1).Value = "Brake_pressure_Front"
$ws.Cells.Item(36, 2).Value = 8
$ws.Cells.Item(36, 3).Value = 8
$ws.Cells.Item(36, 4).Value = "Intel"
$ws.Cells.Item(36, 5).Value = $false
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 0

$ws.Cells.Item(37, 1).Value = "Brake_pressure_Rear"
$ws.Cells.Item(37, 2).Value = 16
$ws.Cells.Item(37, 3).Value = 8
$ws.Cells.Item(37, 4).Value = "Intel"
$ws.Cells.Item(37, 5).Value = $false
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(37, 7).Value = 0
